# Weekly update: insert a new pricing observation (Coliflor, Vega Monumental
# Concepción) dated 2021-09-10 (serial 44449, Región Metropolitana) ahead of
# the existing history, pushing the rest of the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 120:121 - everything that used to live at
# row 120 onward (through the old row 135) shifts down to 122..137.
$ws.Rows("120:121").Insert()

# Row 120 - Primera, 2021-09-10, Región Metropolitana
$ws.Range("A120").Value = 11
$ws.Range("B120").Value = "Vega Monumental Concepción"
$ws.Range("C120").Value = "Bíobío"
$ws.Range("D120").Value = 44449
$ws.Range("E120").Value = 8
$ws.Range("F120").Value = 100112008
$ws.Range("G120").Value = "Coliflor"
$ws.Range("H120").Value = "Sin especificar"
$ws.Range("I120").Value = "Primera"
$ws.Range("J120").Value = 1000
$ws.Range("K120").Value = 700
$ws.Range("L120").Value = 800
$ws.Range("M120").Value = 750
$ws.Range("N120").Value = "$/unidad"
$ws.Range("O120").Value = "Región Metropolitana"
$ws.Range("P120").Value = 750
$ws.Range("Q120").Value = 1
$ws.Range("R120").Value = "Hortaliza"

# Row 121 - Segunda, 2021-09-10, Región Metropolitana
$ws.Range("A121").Value = 11
$ws.Range("B121").Value = "Vega Monumental Concepción"
$ws.Range("C121").Value = "Bíobío"
$ws.Range("D121").Value = 44449
$ws.Range("E121").Value = 8
$ws.Range("F121").Value = 100112008
$ws.Range("G121").Value = "Coliflor"
$ws.Range("H121").Value = "Sin especificar"
$ws.Range("I121").Value = "Segunda"
$ws.Range("J121").Value = 500
$ws.Range("K121").Value = 600
$ws.Range("L121").Value = 600
$ws.Range("M121").Value = 600
$ws.Range("N121").Value = "$/unidad"
$ws.Range("O121").Value = "Región Metropolitana"
$ws.Range("P121").Value = 600
$ws.Range("Q121").Value = 1
$ws.Range("R121").Value = "Hortaliza"
